# Update "nota_view" (column J) for the week 11/09/2022 - 17/09/2022.
# Every student who currently has a nota_view of 5 gets it changed to 4.
# Rows whose nota_view is already 0 stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 10).End(-4162).Row  # xlUp = -4162, column J = 10

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)  # column J
    if ($cell.Value2 -eq 5) {
        $cell.Value2 = 4
    }
}
